$d = $word.ActiveDocument

# The document's paragraphs each carry a <w:contextualSpacing w:val="0"/>
# element inside <w:pPr> (right after <w:shd .../> and before <w:rPr>).
# The target edit removes that element from every paragraph, leaving the
# rest of each paragraph's formatting/content untouched.
#
# ContextualSpacing isn't exposed as a ParagraphFormat property in this
# object model, so we round-trip each paragraph's own WordOpenXML,
# strip the <w:contextualSpacing .../> element from it, and write the
# result back with Range.InsertXML (which replaces that range's
# contents in place).

$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $r = $p.Range
    $xml = $r.WordOpenXML
    if ($xml -match "<w:contextualSpacing[^/]*/>") {
        $newXml = $xml -replace "<w:contextualSpacing[^/]*/>", ""
        $r.InsertXML($newXml)
    }
}
